$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits at the end of the
#    document body (right after the inline picture).
# ---------------------------------------------------------------------
try {
    $goBack = $d.Bookmarks("_GoBack")
    if ($goBack.Start -ne $null) {
        $goBack.Delete()
    }
} catch {
}

# ---------------------------------------------------------------------
# 2) Mint the even/first/default footers + the first-page header by
#    writing into the (currently non-existent) footers.  Word lazily
#    creates footer1.xml/footer2.xml/footer3.xml and header3.xml, wires
#    up the relationships + sectPr header/footerReference elements and
#    the content types the moment any footer story is touched.
# ---------------------------------------------------------------------
$sec = $d.Sections(1)
$primaryFooter = $sec.Footers(1)
$primaryFooter.Range.Text = ""

# ---------------------------------------------------------------------
# 3) Update the printed date from "May 7, 2013" to "June 2, 2013" in
#    both headers.
#    Headers(1) -> default header (header2.xml, simple paragraphs).
#    Headers(3) -> even-page header (header1.xml, table + PAGE field).
# ---------------------------------------------------------------------
$defaultHeader = $sec.Headers(1)
$defaultHeader.Range.Find.Execute("May 7, 2013", $false, $false, $false, `
    $false, $false, $true, 1, $false, "June 2, 2013", 2) | Out-Null

$evenHeader = $sec.Headers(3)

# The even header's Range wraps a table; a direct Find/Replace into the
# table cell is silently dropped unless the header story has already
# been "touched" by a successful Range mutation first.  Insert and
# immediately discard a harmless marker at the end of the header to
# warm it up.
$touch = $evenHeader.Range.Duplicate
$touch.Collapse(0)
$touch.InsertAfter("zzTouchzz")

$evenHeader2 = $sec.Headers(3)
$evenHeader2.Range.Find.Execute("zzTouchzz", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 2) | Out-Null

$evenHeader3 = $sec.Headers(3)
$evenHeader3.Range.Find.Execute("May 7, 2013", $false, $false, $false, `
    $false, $false, $true, 1, $false, "June 2, 2013", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Re-create the "_GoBack" bookmark, this time anchored at the very
#    end of the even-page header (its final, empty paragraph).
# ---------------------------------------------------------------------
$evenHeader4 = $sec.Headers(3)
$mark = $evenHeader4.Range.Duplicate
$mark.Collapse(0)
$mark.InsertAfter("zzGoBackzz")

$evenHeader5 = $sec.Headers(3)
$markRange = $evenHeader5.Range.Duplicate
$markRange.Find.Execute("zzGoBackzz") | Out-Null
$markRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $markRange) | Out-Null

$evenHeader6 = $sec.Headers(3)
$evenHeader6.Range.Find.Execute("zzGoBackzz", $false, $false, $false, `
    $false, $false, $true, 1, $false, "", 2) | Out-Null

"done"
